{"js": "// Replace the raw markdown-style image placeholders (leftover from the\n// source markdown -> docx conversion, e.g. \"](images/foo.png)\") with\n// proper \"[INSERT FIGURE n.n NEAR HERE]\" typesetting notes.\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\n// Exact paragraph-text -> replacement-text map. Each of these paragraphs\n// in the source document consists of a single run whose entire text is\n// the markdown image-link remnant, so a full paragraph text replacement\n// reproduces the change exactly (note the Tesla-coil-1 paragraph loses\n// its leading sentence too, per the target edit).\nconst replacements = {\n  \"](images/me_cover_aerophone.jpg)\": \"[INSERT FIGURE 3.2 NEAR HERE]\",\n  \"](images/electrical_baths.png)\": \"[INSERT FIGURE 3.1 NEAR HERE]\",\n  \"Ever since this magazine was started, subscribers have been clamoring for an article on Tesla coil experiments.](images/tesla_coil1.png)\":\n    \"[INSERT FIGURE 3.3a NEAR HERE]\",\n  \"](images/tesla_coil2.png)\": \"[INSERT FIGURE 3.3b NEAR HERE]\",\n  \"](images/tesla_coil3.png)\": \"[INSERT FIGURE 3.3c NEAR HERE]\",\n  \"](images/tesla_coil4.png)\": \"[INSERT FIGURE 3.3d NEAR HERE]\",\n  \"](images/tesla_coil5.png)\": \"[INSERT FIGURE 3.3e NEAR HERE]\",\n};\n\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  const para = paragraphs.items[i];\n  const newText = replacements[para.text];\n  if (newText !== undefined) {\n    para.insertText(newText, \"Replace\");\n  }\n}\n\nawait context.sync();\n", "ps1": "# Replace the raw markdown-style image placeholders (leftover from the\n# source markdown -> docx conversion, e.g. \"](images/foo.png)\") with\n# proper \"[INSERT FIGURE n.n NEAR HERE]\" typesetting notes.\n$d = $word.ActiveDocument\n\n# Exact paragraph-text -> replacement-text map (keys are the plain\n# paragraph text, without the trailing paragraph mark). Each of these\n# paragraphs consists of a single run whose entire text is the markdown\n# image-link remnant, so a full paragraph text replacement reproduces the\n# change exactly (the Tesla-coil-1 paragraph loses its leading sentence\n# too, per the target edit). Using Paragraph.Range.Text (rather than\n# Find/Replace) keeps the run's xml:space=\"preserve\" markup intact.\n$replacements = @{\n    \"](images/me_cover_aerophone.jpg)\" = \"[INSERT FIGURE 3.2 NEAR HERE]\";\n    \"](images/electrical_baths.png)\" = \"[INSERT FIGURE 3.1 NEAR HERE]\";\n    \"Ever since this magazine was started, subscribers have been clamoring for an article on Tesla coil experiments.](images/tesla_coil1.png)\" = \"[INSERT FIGURE 3.3a NEAR HERE]\";\n    \"](images/tesla_coil2.png)\" = \"[INSERT FIGURE 3.3b NEAR HERE]\";\n    \"](images/tesla_coil3.png)\" = \"[INSERT FIGURE 3.3c NEAR HERE]\";\n    \"](images/tesla_coil4.png)\" = \"[INSERT FIGURE 3.3d NEAR HERE]\";\n    \"](images/tesla_coil5.png)\" = \"[INSERT FIGURE 3.3e NEAR HERE]\"\n}\n\nforeach ($p in $d.Paragraphs) {\n    $paraText = $p.Range.Text.TrimEnd([char]13, [char]7)\n    if ($replacements.ContainsKey($paraText)) {\n        $p.Range.Text = $replacements[$paraText]\n    }\n}\n"}
